# Regenerate save_data: column G ("K") is recalculated (K = strikeouts/knockdowns stat,
# now derived differently -- "use K instead of Strike#") and the newly computed
# s_vals are written back into the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> newly calculated K value (s_val) for each existing data row.
$sVals = @{
    2  = 0
    3  = 2
    4  = 0
    5  = 2
    6  = 0
    7  = 0
    8  = 2
    9  = 1
    10 = 2
    11 = 5
    12 = 0
    13 = 3
    14 = 0
    15 = 0
    16 = 1
    17 = 2
    18 = 2
    19 = 0
    20 = 1
    21 = 0
    22 = 2
    23 = 1
    24 = 1
    25 = 0
    26 = 0
    27 = 3
    28 = 1
    29 = 0
    30 = 1
    31 = 1
    32 = 0
    33 = 2
}

foreach ($row in $sVals.Keys) {
    $ws.Range("G$row").Value = $sVals[$row]
}
